# Insert a new data row at row 117 (pushing the existing row 117 and all
# subsequent rows down by one), then populate the newly inserted row with
# the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 117..194 down to 118..195 by inserting a new row at 117.
$ws.Rows.Item(117).Insert()

# Populate the new row 117 with the new record's values.
$ws.Cells.Item(117, 1).Value2 = 4
$ws.Cells.Item(117, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(117, 3).Value2 = "Los Lagos"
$ws.Cells.Item(117, 4).Value2 = 44603
$ws.Cells.Item(117, 5).Value2 = 10
$ws.Cells.Item(117, 6).Value2 = 100112032
$ws.Cells.Item(117, 7).Value2 = "Zapallo italiano"
$ws.Cells.Item(117, 8).Value2 = "Sin especificar"
$ws.Cells.Item(117, 9).Value2 = "Primera"
$ws.Cells.Item(117, 10).Value2 = 200
$ws.Cells.Item(117, 11).Value2 = 15000
$ws.Cells.Item(117, 12).Value2 = 16000
$ws.Cells.Item(117, 13).Value2 = 15500
$ws.Cells.Item(117, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(117, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(117, 16).Value2 = 310
$ws.Cells.Item(117, 17).Value2 = 50
$ws.Cells.Item(117, 18).Value2 = "Hortaliza"
